$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column Q data by copying formatting from column P, then updating values.
$ws.Range("P2").Copy($ws.Range("Q2"))

$ws.Range("P3").Copy($ws.Range("Q3"))
$ws.Range("Q3").Value = 2023

$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 279.01945525291825

$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 1792.7

$ws.Range("P6").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = 6425

# Adjust row heights for rows 4 and 5 to become custom heights.
$ws.Rows.Item(4).RowHeight = 27
$ws.Rows.Item(5).RowHeight = 27.75

# Reset selection back to the default (A1), matching the committed workbook state.
$ws.Range("A1").Select()
